# Fruta / hortaliza, semanal
# Insert a new weekly record at row 10, pushing the existing rows 10-34
# down to 11-35, then populate the new row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10 (shifts rows 10:34 -> 11:35,
# carrying formatting/styles along, same as Excel's Rows.Insert UI action).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly data point.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 45002
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112044
$ws.Range("G10").Value = "Perejil"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 1200
$ws.Range("L10").Value = 1200
$ws.Range("M10").Value = 1200
$ws.Range("N10").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 1200
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
